$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version: -> 8 (numeric, centered like the other numeric cells such as B2)
$ws.Range("D2").Value = 8
$ws.Range("D2").HorizontalAlignment = -4108

# Descripcion: -> Pruebas Basicas
$ws.Range("D3").Value = "Pruebas Basicas"

# Probado en: -> Google Chrome
$ws.Range("D4").Value = "Google Chrome"

# Update the active selection to D2
$ws.Range("D2").Select() | Out-Null
